$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("G2").Value = 1.67
$ws.Range("H2").Value = 3.5
$ws.Range("I2").Value = 5.75
$ws.Range("J2").Value = 2.38
$ws.Range("K2").Value = 2.05
$ws.Range("L2").Value = 6
$ws.Range("Q2").Value = 2.35
$ws.Range("R2").Value = 1.57
$ws.Range("S2").Value = 1.5
$ws.Range("T2").Value = 2.5
$ws.Range("U2").Value = 2.25
$ws.Range("V2").Value = 1.57
$ws.Range("W2").Value = 5
$ws.Range("X2").Value = 6.5
$ws.Range("Z2").Value = 12
$ws.Range("AA2").Value = 17
$ws.Range("AD2").Value = 7
$ws.Range("AG2").Value = 11
$ws.Range("AH2").Value = 26
$ws.Range("AI2").Value = 19
$ws.Range("AJ2").Value = 67
$ws.Range("AK2").Value = 51
$ws.Range("AN2").Value = 3.4
$ws.Range("AO2").Value = 9
$ws.Range("AP2").Value = 23
$ws.Range("AQ2").Value = 34
$ws.Range("AT2").Value = 2.5
$ws.Range("AV2").Value = 81
$ws.Range("AW2").Value = 7
$ws.Range("AX2").Value = 34
$ws.Range("AZ2").Value = 126
$ws.Range("BA2").Value = 201
$ws.Range("BB2").Value = 451

# Row 4
$ws.Range("G4").Value = 3.3
$ws.Range("H4").Value = 3.3
$ws.Range("I4").Value = 2.25
$ws.Range("J4").Value = 4
$ws.Range("L4").Value = 3
$ws.Range("N4").Value = 8.5
$ws.Range("W4").Value = 9
$ws.Range("Y4").Value = 12
$ws.Range("AA4").Value = 29
$ws.Range("AB4").Value = 41
$ws.Range("AG4").Value = 7
$ws.Range("AH4").Value = 10
$ws.Range("AJ4").Value = 21
$ws.Range("AK4").Value = 19
$ws.Range("AL4").Value = 29
$ws.Range("AO4").Value = 19
$ws.Range("AQ4").Value = 67

# Row 5
$ws.Range("G5").Value = 3.5
$ws.Range("I5").Value = 2.2
$ws.Range("J5").Value = 4.33
$ws.Range("L5").Value = 3.1
$ws.Range("Q5").Value = 2.7
$ws.Range("R5").Value = 1.44
$ws.Range("X5").Value = 15
$ws.Range("AJ5").Value = 21
$ws.Range("AU5").Value = 9.5
$ws.Range("AZ5").Value = 51
$ws.Range("BB5").Value = 301

# Row 8
$ws.Range("V8").Value = 1.57
